$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction11")

# Clear out the old wide row (A1:P1) then write the new, smaller values.
$ws.Range("A1:P1").Clear()

$ws.Range("A1").Value = 22
$ws.Range("B1").Value = 23
